$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.8200883333333334
$ws.Range("H2").Value = 2.460265
$ws.Range("I2").Value = 0.2405117342909232
$ws.Range("J2").Value = 0.2405117342909232
$ws.Range("Q2").Value = 0.01671558713555556
$ws.Range("R2").Value = 0.15044028422
$ws.Range("S2").Value = 0.2405117342909232
$ws.Range("T2").Value = 0.2405117342909232

# Row 3 updates
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.589676
$ws.Range("H3").Value = 7.769028
$ws.Range("I3").Value = 0.7594882657090768
$ws.Range("J3").Value = 0.7594882657090768
$ws.Range("Q3").Value = 0.05278450268266666
$ws.Range("R3").Value = 0.475060524144
$ws.Range("S3").Value = 0.7594882657090768
$ws.Range("T3").Value = 0.7594882657090768
